$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Version 11.03.10, 2015-06-29",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Version 11.04.03, 2015-07-15",
    2
)
